$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.179.11'
$ws.Range("E2").Value = '  +5.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.464.43'
$ws.Range("E3").Value = '  +6.47%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.10'
$ws.Range("E5").Value = '  +4.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.56'
$ws.Range("E6").Value = '  +10.68%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.592'
$ws.Range("E8").Value = '  +2.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.462.54'
$ws.Range("E9").Value = '  +6.34%  '
$ws.Range("E10").Value = '  +5.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.72'
$ws.Range("E11").Value = '  +2.72%  '
$ws.Range("E12").Value = '  +1.31%  '
$ws.Range("E13").Value = '  +5.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.55'
$ws.Range("E14").Value = '  +13.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.906.24'
$ws.Range("E15").Value = '  +6.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.087.09'
$ws.Range("E16").Value = '  +4.99%  '
$ws.Range("E17").Value = '  +6.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.465.88'
$ws.Range("E18").Value = '  +6.98%  '
$ws.Range("E19").Value = '  +6.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.70'
$ws.Range("E20").Value = '  +9.14%  '
$ws.Range("E21").Value = '  +5.57%  '
$ws.Range("E22").Value = '  +3.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("E24").Value = '  +2.39%  '
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.51'
$ws.Range("E27").Value = '  +10.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.16'
$ws.Range("E28").Value = '  +4.49%  '
$ws.Range("E29").Value = '  +7.95%  '
$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.84'
$ws.Range("E30").Value = '  +14.41%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0812'
$ws.Range("E31").Value = '  +12.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.84'
$ws.Range("E32").Value = '  +7.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '175.77'
$ws.Range("E33").Value = '  +2.61%  '
$ws.Range("E34").Value = '  +12.51%  '
$ws.Range("E35").Value = '  +4.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.92'
$ws.Range("E36").Value = '  +5.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '369.89'
$ws.Range("E37").Value = '  +17.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.46'
$ws.Range("E38").Value = '  +9.82%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("E41").Value = '  +12.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.42'
$ws.Range("E42").Value = '  +6.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '151.03'
$ws.Range("E43").Value = '  +9.92%  '
$ws.Range("E44").Value = '  +6.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.62'
$ws.Range("E45").Value = '  +8.44%  '
$ws.Range("E46").Value = '  +6.03%  '
$ws.Range("E47").Value = '  +2.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0520'
$ws.Range("E48").Value = '  +5.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0244'
$ws.Range("E49").Value = '  +12.20%  '
$ws.Range("E50").Value = '  +5.42%  '
$ws.Range("E51").Value = '  +7.02%  '
